$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")
$ws.Range("F6").Value = "Leandro"
$ws.Range("E6").Value = "Desenvolvendo "
$ws.Range("G16").Orientation = 0
$ws.Range("G16").Select()
